$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 (ano = 2025) metrics with refreshed figures
$ws.Range("C6").Value = 344
$ws.Range("D6").Value = 280
$ws.Range("E6").Value = 64
$ws.Range("F6").Value = 62.36080178173719
$ws.Range("G6").Value = 18.6046511627907
$ws.Range("H6").Value = 81.3953488372093
